$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Update cell values (dates + rolling-window numeric shifts) ---
$ws.Range("N3").Value2 = 46023
$ws.Range("N4").Value2 = 46023
$ws.Range("N6").Value2 = 46023
$ws.Range("N7").Value2 = 46023
$ws.Range("N8").Value2 = 46023
$ws.Range("N9").Value2 = 46023
$ws.Range("N15").Value2 = 46023
$ws.Range("C28").Value2 = 45992
$ws.Range("F28").Value2 = -0.01425751734772551
$ws.Range("G28").Value2 = 0.05423155704392491
$ws.Range("H28").Value2 = -0.02093605859677161
$ws.Range("I28").Value2 = 0.006436255758670795
$ws.Range("J28").Value2 = 0.03004963172206243
$ws.Range("C29").Value2 = 45992
$ws.Range("F29").Value2 = 0.1000416438657897
$ws.Range("G29").Value2 = 0.1245404730526462
$ws.Range("H29").Value2 = 0.04877483240471108
$ws.Range("I29").Value2 = 0.07412067603746038
$ws.Range("J29").Value2 = 0.07661265288383932
$ws.Range("N29").Value2 = 46070
$ws.Range("Q29").Value2 = 2.13
$ws.Range("R29").Value2 = 2.12
$ws.Range("S29").Value2 = 2.13
$ws.Range("T29").Value2 = 2.15
$ws.Range("U29").Value2 = 2.17
$ws.Range("C30").Value2 = 45992
$ws.Range("F30").Value2 = -0.02458405517602069
$ws.Range("G30").Value2 = 0.06576508785864577
$ws.Range("H30").Value2 = -0.0128376635658648
$ws.Range("I30").Value2 = 0.001174064535676367
$ws.Range("J30").Value2 = 0.01907672443132968
$ws.Range("N30").Value2 = 46070
$ws.Range("Q30").Value2 = 2.26
$ws.Range("R30").Value2 = 2.27
$ws.Range("S30").Value2 = 2.29
$ws.Range("U30").Value2 = 2.32
$ws.Range("C31").Value2 = 45992
$ws.Range("F31").Value2 = 0.09489005566012468
$ws.Range("G31").Value2 = 0.1263576451529767
$ws.Range("H31").Value2 = 0.05011084527755218
$ws.Range("I31").Value2 = 0.06502168244015354
$ws.Range("J31").Value2 = 0.06671073894520346
$ws.Range("C32").Value2 = 46023
$ws.Range("F32").Value2 = 0.007001897085101128
$ws.Range("G32").Value2 = 0.002486745086434317
$ws.Range("H32").Value2 = 0.001190041182338009
$ws.Range("I32").Value2 = -0.004413706579460941
$ws.Range("J32").Value2 = 0.0007990183488855163
$ws.Range("C33").Value2 = 46023
$ws.Range("F33").Value2 = 0.02275028056847218
$ws.Range("G33").Value2 = 0.01298051477514099
$ws.Range("H33").Value2 = 0.02099856484628742
$ws.Range("I33").Value2 = 0.01797033261452011
$ws.Range("J33").Value2 = 0.01901142589200902
$ws.Range("N33").Value2 = 46023
$ws.Range("C34").Value2 = 46023
$ws.Range("F34").Value2 = 76.2119
$ws.Range("G34").Value2 = 75.7433
$ws.Range("H34").Value2 = 75.64619999999999
$ws.Range("I34").Value2 = 75.6474
$ws.Range("J34").Value2 = 76.0745
$ws.Range("C36").Value2 = 45992
$ws.Range("F36").Value2 = 1404
$ws.Range("G36").Value2 = 1322
$ws.Range("H36").Value2 = 1272
$ws.Range("I36").Value2 = 1328
$ws.Range("J36").Value2 = 1291
$ws.Range("C37").Value2 = 45992
$ws.Range("F37").Value2 = -0.0726552179656539
$ws.Range("G37").Value2 = 0.02084942084942085
$ws.Range("H37").Value2 = -0.05917159763313609
$ws.Range("I37").Value2 = -0.02137067059690494
$ws.Range("J37").Value2 = -0.07189072609633357
$ws.Range("C38").Value2 = 45992
$ws.Range("F38").Value2 = 1448
$ws.Range("G38").Value2 = 1388
$ws.Range("H38").Value2 = 1411
$ws.Range("I38").Value2 = 1415
$ws.Range("J38").Value2 = 1330
$ws.Range("C39").Value2 = 45992
$ws.Range("F39").Value2 = -0.02162162162162162
$ws.Range("G39").Value2 = -0.07957559681697612
$ws.Range("H39").Value2 = -0.0119047619047619
$ws.Range("I39").Value2 = -0.01324965132496513
$ws.Range("J39").Value2 = -0.0989159891598916
$ws.Range("N47").Value2 = 46069
$ws.Range("N48").Value2 = 46066
$ws.Range("Q48").Value2 = 3.4
$ws.Range("R48").Value2 = 3.47
$ws.Range("S48").Value2 = 3.52
$ws.Range("T48").Value2 = 3.45
$ws.Range("U48").Value2 = 3.48
$ws.Range("N49").Value2 = 46066
$ws.Range("Q49").Value2 = 3.61
$ws.Range("R49").Value2 = 3.67
$ws.Range("S49").Value2 = 3.75
$ws.Range("T49").Value2 = 3.7
$ws.Range("U49").Value2 = 3.75
$ws.Range("N50").Value2 = 46066
$ws.Range("Q50").Value2 = 4.04
$ws.Range("R50").Value2 = 4.09
$ws.Range("S50").Value2 = 4.18
$ws.Range("T50").Value2 = 4.16
$ws.Range("N52").Value2 = 46066
$ws.Range("Q52").Value2 = 5.76
$ws.Range("R52").Value2 = 5.77
$ws.Range("S52").Value2 = 5.85
$ws.Range("T52").Value2 = 5.82
$ws.Range("U52").Value2 = 5.86

# --- Step 2: Re-point "latest updated" highlight (yellow fill) ---
# Remove yellow highlight from N-column cells that are no longer the newest period
$ws.Range("C3").Copy()  # C3 has the plain (non-highlighted) date style
$ws.Range("N3").PasteSpecial(-4122)
$ws.Range("N4").PasteSpecial(-4122)
$ws.Range("N6").PasteSpecial(-4122)
$ws.Range("N7").PasteSpecial(-4122)
$ws.Range("N8").PasteSpecial(-4122)
$ws.Range("N9").PasteSpecial(-4122)
$ws.Range("N15").PasteSpecial(-4122)
$ws.Range("N33").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Add yellow highlight to C-column cells that were just refreshed
$ws.Range("N29").Copy()  # N29 has the highlighted (yellow) date style
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C29").PasteSpecial(-4122)
$ws.Range("C30").PasteSpecial(-4122)
$ws.Range("C31").PasteSpecial(-4122)
$ws.Range("C32").PasteSpecial(-4122)
$ws.Range("C33").PasteSpecial(-4122)
$ws.Range("C34").PasteSpecial(-4122)
$ws.Range("C36").PasteSpecial(-4122)
$ws.Range("C37").PasteSpecial(-4122)
$ws.Range("C38").PasteSpecial(-4122)
$ws.Range("C39").PasteSpecial(-4122)
$excel.CutCopyMode = $false
